# Auto-generated Excel COM-interop script to refresh Leve profit market-data columns (H:N)
# across multiple worksheets, matching the scheduled market-data runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 126 (sheet ALC)
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
# Row 138 (sheet ALC)
$ws.Range("H138").Value = 2366.3696
$ws.Range("I138").Value = 2742.5715
$ws.Range("J138").Value = 2298.8462
$ws.Range("K138").Value = 8227.7145
$ws.Range("L138").Value = 6896.5386
$ws.Range("M138").Value = -3087.7145
$ws.Range("N138").Value = -17176.5386

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (sheet ARM)
$ws.Range("H2").Value = 1851.0834
$ws.Range("I2").Value = 1851.0834
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1851.0834
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1738.0834
$ws.Range("N2").ClearContents()
# Row 97 (sheet ARM)
$ws.Range("H97").Value = 902
$ws.Range("I97").Value = 902
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 902
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -406
$ws.Range("N97").ClearContents()
# Row 110 (sheet ARM)
$ws.Range("H110").Value = 1360.6923
$ws.Range("I110").Value = 1390.75
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1390.75
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 654.25
$ws.Range("N110").Value = -5090
# Row 112 (sheet ARM)
$ws.Range("H112").Value = 79387
$ws.Range("J112").Value = 79387
$ws.Range("L112").Value = 79387
$ws.Range("N112").Value = -82341
# Row 116 (sheet ARM)
$ws.Range("H116").Value = 1851.0834
$ws.Range("I116").Value = 1851.0834
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1851.0834
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 442.9166
$ws.Range("N116").ClearContents()
# Row 132 (sheet ARM)
$ws.Range("H132").Value = 4088.4
$ws.Range("I132").Value = 5512.4
$ws.Range("K132").Value = 16537.2
$ws.Range("M132").Value = -14007.2

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (sheet BSM)
$ws.Range("H3").Value = 1851.0834
$ws.Range("I3").Value = 1851.0834
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1851.0834
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1737.0834
$ws.Range("N3").ClearContents()
# Row 99 (sheet BSM)
$ws.Range("H99").Value = 1173.9333
$ws.Range("I99").Value = 1025.75
$ws.Range("J99").Value = 1766.6666
$ws.Range("K99").Value = 1025.75
$ws.Range("L99").Value = 1766.6666
$ws.Range("M99").Value = 472.25
$ws.Range("N99").Value = -4762.6666
# Row 105 (sheet BSM)
$ws.Range("H105").Value = 7815135.5
$ws.Range("I105").Value = 8931069
$ws.Range("K105").Value = 8931069
$ws.Range("M105").Value = -8929322
# Row 110 (sheet BSM)
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 134 (sheet BSM)
$ws.Range("H134").Value = 2362.122
$ws.Range("I134").Value = 2128.2415
$ws.Range("K134").Value = 6384.7245
$ws.Range("M134").Value = -3849.7245

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (sheet CRP)
$ws.Range("H4").Value = 27545.637
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 27545.637
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 27545.637
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -27769.637
# Row 68 (sheet CRP)
$ws.Range("H68").Value = 23498.615
$ws.Range("J68").Value = 23498.615
$ws.Range("L68").Value = 23498.615
$ws.Range("N68").Value = -24996.615
# Row 70 (sheet CRP)
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 71 (sheet CRP)
$ws.Range("H71").Value = 23498.615
$ws.Range("J71").Value = 23498.615
$ws.Range("L71").Value = 70495.845
$ws.Range("N71").Value = -77983.845
# Row 73 (sheet CRP)
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 75 (sheet CRP)
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41996
# Row 78 (sheet CRP)
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129984
# Row 88 (sheet CRP)
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91 (sheet CRP)
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 105 (sheet CRP)
$ws.Range("H105").Value = 1669.9231
$ws.Range("I105").Value = 1601
$ws.Range("J105").Value = 1825
$ws.Range("K105").Value = 1601
$ws.Range("L105").Value = 1825
$ws.Range("M105").Value = 146
$ws.Range("N105").Value = -5319

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (sheet CUL)
$ws.Range("H5").Value = 535
$ws.Range("I5").Value = 476.875
$ws.Range("K5").Value = 1430.625
$ws.Range("M5").Value = -1318.625
# Row 11 (sheet CUL)
$ws.Range("H11").Value = 166.66667
$ws.Range("I11").Value = 166.66667
$ws.Range("K11").Value = 500.00001
$ws.Range("M11").Value = -360.00001
# Row 113 (sheet CUL)
$ws.Range("H113").Value = 1427.75
$ws.Range("I113").Value = 645.75
$ws.Range("J113").Value = 1818.75
$ws.Range("K113").Value = 1937.25
$ws.Range("L113").Value = 5456.25
$ws.Range("M113").Value = 232.75
$ws.Range("N113").Value = -9796.25
# Row 131 (sheet CUL)
$ws.Range("H131").Value = 989.9429
$ws.Range("I131").Value = 376.66666
$ws.Range("J131").Value = 1047.4375
$ws.Range("K131").Value = 1129.99998
$ws.Range("L131").Value = 3142.3125
$ws.Range("M131").Value = 3910.00002
$ws.Range("N131").Value = -13222.3125
# Row 135 (sheet CUL)
$ws.Range("H135").Value = 535
$ws.Range("I135").Value = 476.875
$ws.Range("K135").Value = 4291.875
$ws.Range("M135").Value = -1756.875

$ws = $wb.Worksheets.Item("GSM")
# Row 111 (sheet GSM)
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# Row 126 (sheet GSM)
$ws.Range("H126").Value = 2111.9443
$ws.Range("I126").Value = 1900.0769
$ws.Range("J126").Value = 2662.8
$ws.Range("K126").Value = 5700.2307
$ws.Range("L126").Value = 7988.400000000001
$ws.Range("M126").Value = -3230.2307
$ws.Range("N126").Value = -12928.4

$ws = $wb.Worksheets.Item("LTW")
# Row 110 (sheet LTW)
$ws.Range("H110").Value = 10000
$ws.Range("J110").Value = 10000
$ws.Range("L110").Value = 10000
$ws.Range("N110").Value = -18180

$ws = $wb.Worksheets.Item("WVR")
# Row 17 (sheet WVR)
$ws.Range("H17").Value = 28182.8
$ws.Range("I17").Value = 301.33334
$ws.Range("K17").Value = 301.33334
$ws.Range("M17").Value = -129.33334
# Row 40 (sheet WVR)
$ws.Range("H40").Value = 15000
$ws.Range("J40").Value = 15000
$ws.Range("L40").Value = 15000
$ws.Range("N40").Value = -15298
# Row 113 (sheet WVR)
$ws.Range("H113").Value = 1128.875
$ws.Range("I113").Value = 1128.875
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3386.625
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1216.625
$ws.Range("N113").ClearContents()
